$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 290
$ws1.Range("F4").Value = 2759
$ws1.Range("F6").Value = 581

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 290
$ws4.Range("F6").Value = 2759
$ws4.Range("F8").Value = 581
